# Update the "Models" sheet with the custom mesh function results for the
# last species (Market squid) and refresh the best-model picks per size
# class, then move the active selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")

# Best model (by size class) per species - refreshed values using the
# custom mesh function output and new scripts for the last species
# (Market squid).
$ws.Range("B4").Value = "uvint 100m*"
$ws.Range("B7").Value = "uvint 100m"

# New header label over the notes column.
$ws.Range("E1").Value = "base: no variable"
$ws.Range("E1").Font.Bold = $true

$ws.Range("B2").Value = "u_vint_50m"
$ws.Range("B3").Value = "u_vint_100m"

$ws.Range("C2").Value = "vmax_cu"
$ws.Range("C3").Value = "vmax_cu"

$ws.Range("C4").Value = "spice*"
$ws.Range("C4").Style = "Normal"

$ws.Range("B5").Value = "ios26"
$ws.Range("C5").Value = "spice*"

$ws.Range("B6").Value = "v_cu"
$ws.Range("C6").Value = "v_cu"

$ws.Range("C7").Value = "base"

# Move the active cell/selection to C3.
$ws.Range("C3").Select()
